$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-04-23 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-24 Wednesday", 2) | Out-Null

# Update the multiplication problems in the single table, addressed by
# (row, column) so that duplicate expressions are replaced unambiguously
# and in document order.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "33×17="
$t.Cell(1, 2).Range.Text = "33×61="
$t.Cell(1, 3).Range.Text = "99×42="
$t.Cell(1, 4).Range.Text = "25×14="
$t.Cell(1, 5).Range.Text = "27×35="
$t.Cell(5, 1).Range.Text = "13×88="
$t.Cell(5, 2).Range.Text = "79×74="
$t.Cell(5, 3).Range.Text = "68×36="
$t.Cell(5, 4).Range.Text = "49×29="
$t.Cell(5, 5).Range.Text = "47×47="
$t.Cell(10, 1).Range.Text = "72×99="
$t.Cell(10, 2).Range.Text = "58×77="
$t.Cell(10, 3).Range.Text = "96×74="
$t.Cell(10, 4).Range.Text = "91×67="
$t.Cell(10, 5).Range.Text = "29×58="
$t.Cell(15, 1).Range.Text = "57×28="
$t.Cell(15, 2).Range.Text = "62×28="
$t.Cell(15, 3).Range.Text = "84×25="
$t.Cell(15, 4).Range.Text = "39×90="
$t.Cell(15, 5).Range.Text = "97×38="
$t.Cell(20, 1).Range.Text = "99×17="
$t.Cell(20, 2).Range.Text = "82×84="
$t.Cell(20, 3).Range.Text = "86×17="
$t.Cell(20, 4).Range.Text = "42×96="
$t.Cell(20, 5).Range.Text = "19×36="
